# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for existing rows 2-29
# from 45586 (2024-10-21) to 45587 (2024-10-22).
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 3).Value = 45587
}

# Row 29 gains an explicit row height (matches default 15, but written out
# explicitly as customHeight).
$ws.Rows.Item(29).RowHeight = 15

# Add new row 30: A 46082-2024
$ws.Cells.Item(30, 1).Value = "A 46082-2024"
$ws.Cells.Item(30, 2).Value = 45581
$ws.Cells.Item(30, 3).Value = 45587
$ws.Cells.Item(30, 4).Value = "OKÄNT"
$ws.Cells.Item(30, 5).Value = "OKÄNT"
$ws.Cells.Item(30, 7).Value = 1
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(30, 14).Value = 0
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 16).Value = 0
$ws.Cells.Item(30, 17).Value = 0
$ws.Cells.Item(30, 18).WrapText = $true
$ws.Rows.Item(30).RowHeight = 15

# Add new row 31: A 46085-2024
$ws.Cells.Item(31, 1).Value = "A 46085-2024"
$ws.Cells.Item(31, 2).Value = 45581
$ws.Cells.Item(31, 3).Value = 45587
$ws.Cells.Item(31, 4).Value = "OKÄNT"
$ws.Cells.Item(31, 5).Value = "OKÄNT"
$ws.Cells.Item(31, 7).Value = 1.9
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(31, 14).Value = 0
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(31, 16).Value = 0
$ws.Cells.Item(31, 17).Value = 0
$ws.Cells.Item(31, 18).WrapText = $true

# Apply the date display format and wrap-text formatting consistent with
# the rest of the table for the new rows.
$ws.Range("B30:C31").NumberFormat = "YYYY-MM-DD"
$ws.Range("R30:R31").WrapText = $true
